$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.213.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.69%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.810.64"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.98%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4317"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.86%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3660"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.83%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.80"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07625"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.141"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.004"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.298"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.449"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.814.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.64"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001076"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06416"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9998"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.215"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.251.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.124"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.022.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.260"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "130.68"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.181"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.967"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09093"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.556"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.87"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02390"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.193"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2163"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6553"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06173"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.200"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.029"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.422"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9990"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6073"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.727"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.70"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.005"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.171"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06995"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.06%  "
